# Refresh cryptos list figures (Price / Volume(1h)) - GitHub Actions scheduled update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.500.13'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '''1.841.85'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''262.49'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').Value = '''0.5303'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('D8').Value = '''0.3087'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.96%  '
$ws.Range('D9').Value = '''0.06903'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.85%  '
$ws.Range('D10').Value = '''18.54'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.17%  '
$ws.Range('D11').Value = '''0.07831'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('D12').Value = '''0.7604'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('D13').Value = '''1.868.16'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = '''89.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.46%  '
$ws.Range('D15').Value = '''5.033'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').Value = '''14.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.51%  '
$ws.Range('D18').Value = '''0.000007945'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '''26.526.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = '''4.627'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('D22').Value = '''6.006'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D23').Value = '''9.311'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.71%  '
$ws.Range('D24').Value = '''142.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').Value = '''2.191'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').Value = '''1.681'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Value = '''17.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').Value = '''111.44'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').Value = '''0.08796'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').Value = '''4.091'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.66%  '
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('D34').Value = '''0.7307'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.34%  '
$ws.Range('D35').Value = '''1.132'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('D36').Value = '''3.101'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +5.96%  '
$ws.Range('D38').Value = '''0.01721'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.78%  '
$ws.Range('D39').Value = '''0.4798'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.02%  '
$ws.Range('D40').Value = '''0.9022'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.47%  '
$ws.Range('D41').Value = '''108.23'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.69%  '
$ws.Range('D42').Value = '''5.885'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.93%  '
$ws.Range('E43').Value = '  -0.08%  '
$ws.Range('D44').Value = '''7.593'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('D45').Value = '''0.4152'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.01%  '
$ws.Range('D46').Value = '''9.052'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').Value = '''0.1240'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.74%  '
$ws.Range('D48').Value = '''34.94'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('D49').Value = '''0.8999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.76%  '
$ws.Range('D50').Value = '''0.05802'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E51').Value = '  +0.80%  '
